$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05008433333333333
$ws.Range("H2").Value = 0.150253
$ws.Range("I2").Value = 0.3054767171413236
$ws.Range("J2").Value = 0.3054767171413236
$ws.Range("M2").Value = 0.5623183333333334
$ws.Range("N2").Value = 1.686955
$ws.Range("O2").Value = 0.2111849380937466
$ws.Range("P2").Value = 0.2111849380937466
$ws.Range("Q2").Value = 0.02816333884611111
$ws.Range("R2").Value = 0.253470049615
$ws.Range("S2").Value = 0.06451208159857139
$ws.Range("T2").Value = 0.06451208159857137

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05008433333333333
$ws.Range("H3").Value = 0.150253
$ws.Range("I3").Value = 0.3054767171413236
$ws.Range("J3").Value = 0.3054767171413236
$ws.Range("O3").Value = 0.3721759989864856
$ws.Range("P3").Value = 0.3721759989864856
$ws.Range("Q3").Value = 0.04963288984744445
$ws.Range("R3").Value = 0.446696008627
$ws.Range("S3").Value = 0.1136911023691842
$ws.Range("T3").Value = 0.1136911023691842

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05008433333333333
$ws.Range("H4").Value = 0.150253
$ws.Range("I4").Value = 0.3054767171413236
$ws.Range("J4").Value = 0.3054767171413236
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.248985
$ws.Range("N4").Value = 0.746955
$ws.Range("O4").Value = 0.09350910097413062
$ws.Range("P4").Value = 0.09350910097413061
$ws.Range("Q4").Value = 0.012470247735
$ws.Range("R4").Value = 0.112232229615
$ws.Range("S4").Value = 0.02856485318841397
$ws.Range("T4").Value = 0.02856485318841397

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05008433333333333
$ws.Range("H5").Value = 0.150253
$ws.Range("I5").Value = 0.3054767171413236
$ws.Range("J5").Value = 0.3054767171413236
$ws.Range("M5").Value = 0.484378
$ws.Range("N5").Value = 1.453134
$ws.Range("O5").Value = 0.1819135743584852
$ws.Range("P5").Value = 0.1819135743584852
$ws.Range("Q5").Value = 0.02425974921133333
$ws.Range("R5").Value = 0.218337742902
$ws.Range("S5").Value = 0.05557036149847413
$ws.Range("T5").Value = 0.05557036149847412

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05008433333333333
$ws.Range("H6").Value = 0.150253
$ws.Range("I6").Value = 0.3054767171413236
$ws.Range("J6").Value = 0.3054767171413236
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2069713333333333
$ws.Range("N6").Value = 0.620914
$ws.Range("O6").Value = 0.07773039864818006
$ws.Range("P6").Value = 0.07773039864818004
$ws.Range("Q6").Value = 0.01036602124911111
$ws.Range("R6").Value = 0.09329419124199999
$ws.Range("S6").Value = 0.02374482700113242
$ws.Range("T6").Value = 0.02374482700113242

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05008433333333333
$ws.Range("H7").Value = 0.150253
$ws.Range("I7").Value = 0.3054767171413236
$ws.Range("J7").Value = 0.3054767171413236
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.169043
$ws.Range("N7").Value = 0.5071289999999999
$ws.Range("O7").Value = 0.06348598893897206
$ws.Range("P7").Value = 0.06348598893897206
$ws.Range("Q7").Value = 0.008466405959666665
$ws.Range("R7").Value = 0.07619765363699998
$ws.Range("S7").Value = 0.01939349148554757
$ws.Range("T7").Value = 0.01939349148554757

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1138703333333333
$ws.Range("H8").Value = 0.341611
$ws.Range("I8").Value = 0.6945232828586764
$ws.Range("J8").Value = 0.6945232828586764
$ws.Range("M8").Value = 0.5623183333333334
$ws.Range("N8").Value = 1.686955
$ws.Range("O8").Value = 0.2111849380937466
$ws.Range("P8").Value = 0.2111849380937466
$ws.Range("Q8").Value = 0.06403137605611112
$ws.Range("R8").Value = 0.576282384505
$ws.Range("S8").Value = 0.1466728564951753
$ws.Range("T8").Value = 0.1466728564951752

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1138703333333333
$ws.Range("H9").Value = 0.341611
$ws.Range("I9").Value = 0.6945232828586764
$ws.Range("J9").Value = 0.6945232828586764
$ws.Range("O9").Value = 0.3721759989864856
$ws.Range("P9").Value = 0.3721759989864856
$ws.Range("Q9").Value = 0.1128439441054444
$ws.Range("R9").Value = 1.015595496949
$ws.Range("S9").Value = 0.2584848966173014
$ws.Range("T9").Value = 0.2584848966173014

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1138703333333333
$ws.Range("H10").Value = 0.341611
$ws.Range("I10").Value = 0.6945232828586764
$ws.Range("J10").Value = 0.6945232828586764
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.248985
$ws.Range("N10").Value = 0.746955
$ws.Range("O10").Value = 0.09350910097413062
$ws.Range("P10").Value = 0.09350910097413061
$ws.Range("Q10").Value = 0.028352004945
$ws.Range("R10").Value = 0.255168044505
$ws.Range("S10").Value = 0.06494424778571666
$ws.Range("T10").Value = 0.06494424778571664

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1138703333333333
$ws.Range("H11").Value = 0.341611
$ws.Range("I11").Value = 0.6945232828586764
$ws.Range("J11").Value = 0.6945232828586764
$ws.Range("M11").Value = 0.484378
$ws.Range("N11").Value = 1.453134
$ws.Range("O11").Value = 0.1819135743584852
$ws.Range("P11").Value = 0.1819135743584852
$ws.Range("Q11").Value = 0.05515628431933333
$ws.Range("R11").Value = 0.496406558874
$ws.Range("S11").Value = 0.1263432128600111
$ws.Range("T11").Value = 0.1263432128600111

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1138703333333333
$ws.Range("H12").Value = 0.341611
$ws.Range("I12").Value = 0.6945232828586764
$ws.Range("J12").Value = 0.6945232828586764
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2069713333333333
$ws.Range("N12").Value = 0.620914
$ws.Range("O12").Value = 0.07773039864818006
$ws.Range("P12").Value = 0.07773039864818004
$ws.Range("Q12").Value = 0.02356789471711111
$ws.Range("R12").Value = 0.212111052454
$ws.Range("S12").Value = 0.05398557164704763
$ws.Range("T12").Value = 0.05398557164704763

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1138703333333333
$ws.Range("H13").Value = 0.341611
$ws.Range("I13").Value = 0.6945232828586764
$ws.Range("J13").Value = 0.6945232828586764
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.169043
$ws.Range("N13").Value = 0.5071289999999999
$ws.Range("O13").Value = 0.06348598893897206
$ws.Range("P13").Value = 0.06348598893897206
$ws.Range("Q13").Value = 0.01924898275766666
$ws.Range("R13").Value = 0.173240844819
$ws.Range("S13").Value = 0.04409249745342449
$ws.Range("T13").Value = 0.04409249745342449
